# Natmi following Dr Hou advice
# The "Sending cluster" dimension now also includes "ECs" (in addition to
# "FAPs" and "sCs"), so the Epha4/Efnb1 LR-pair table grows from 6 data
# rows (2 senders x 3 targets) to 9 data rows (3 senders x 3 targets), and
# every numeric column is recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#         E..T the 16 numeric metric columns)
$data = @{
    2  = @("ECs",  "Epha4", "Efnb1", "ECs",  2, 0.6666666666666666, 2.856403666666667, 8.569210999999999, 0.235832554697756, 0.235832554697756, 3, 1, 3.896501666666667, 11.689505, 0.401720501899026, 0.401720501899026, 11.12998164783945, 100.169834830555, 0.09473877223731204, 0.09473877223731203)
    3  = @("ECs",  "Epha4", "Efnb1", "FAPs", 2, 0.6666666666666666, 2.856403666666667, 8.569210999999999, 0.235832554697756, 0.235832554697756, 3, 1, 3.23724, 9.71172, 0.3337521163387849, 0.3337521163387849, 9.24686420588, 83.22177785291998, 0.07870961423195832, 0.07870961423195832)
    4  = @("ECs",  "Epha4", "Efnb1", "sCs",  2, 0.6666666666666666, 2.856403666666667, 8.569210999999999, 0.235832554697756, 0.235832554697756, 3, 1, 2.565792333333333, 7.697377, 0.2645273817621892, 0.2645273817621892, 7.328938628838555, 65.960447659547, 0.06238416822848566, 0.06238416822848565)
    5  = @("FAPs", "Epha4", "Efnb1", "ECs",  3, 1, 7.000300666666668, 21.000902, 0.5779641054021444, 0.5779641054021444, 3, 1, 3.896501666666667, 11.689505, 0.401720501899026, 0.401720501899026, 27.27668321483445, 245.4901489335101, 0.232180030501771, 0.232180030501771)
    6  = @("FAPs", "Epha4", "Efnb1", "FAPs", 3, 1, 7.000300666666668, 21.000902, 0.5779641054021444, 0.5779641054021444, 3, 1, 3.23724, 9.71172, 0.3337521163387849, 0.3337521163387849, 22.66165333016, 203.95487997144, 0.1928967433458182, 0.1928967433458182)
    7  = @("FAPs", "Epha4", "Efnb1", "sCs",  3, 1, 7.000300666666668, 21.000902, 0.5779641054021444, 0.5779641054021444, 3, 1, 2.565792333333333, 7.697377, 0.2645273817621892, 0.2645273817621892, 17.96131778156156, 161.651860034054, 0.1528873315545552, 0.1528873315545552)
    8  = @("sCs",  "Epha4", "Efnb1", "ECs",  3, 1, 2.255294666666666, 6.765884, 0.1862033399000996, 0.1862033399000996, 3, 1, 3.896501666666667, 11.689505, 0.401720501899026, 0.401720501899026, 8.78775942749111, 79.08983484742001, 0.07480169915994293, 0.07480169915994293)
    9  = @("sCs",  "Epha4", "Efnb1", "FAPs", 3, 1, 2.255294666666666, 6.765884, 0.1862033399000996, 0.1862033399000996, 3, 1, 3.23724, 9.71172, 0.3337521163387849, 0.3337521163387849, 7.300930106719999, 65.70837096048, 0.06214575876100835, 0.06214575876100835)
    10 = @("sCs",  "Epha4", "Efnb1", "sCs",  3, 1, 2.255294666666666, 6.765884, 0.1862033399000996, 0.1862033399000996, 3, 1, 2.565792333333333, 7.697377, 0.2645273817621892, 0.2645273817621892, 5.786617765140889, 52.079559886268, 0.04925588197914831, 0.04925588197914831)
}

# Rows 8-10 are brand new ("ECs" as a Sending cluster is a new category),
# rows 2-7 are the previously existing sender/target combinations with
# their metrics recomputed now that "ECs" is part of the sender pool.
for ($r = 2; $r -le 10; $r++) {
    $row = $data[$r]
    for ($i = 0; $i -lt $row.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $row[$i]
    }
}
